$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.476.62"
$ws.Range("E2").Value = "  -0.45%  "
$ws.Range("D3").Value = "1.571.97"
$ws.Range("E3").Value = "  -1.54%  "
$ws.Range("E4").Value = "  -0.10%  "
$cell = $ws.Range("D5")
$cell.Value = "'207.56"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -0.32%  "
$cell = $ws.Range("D6")
$cell.Value = "'0.496"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -1.28%  "
$ws.Range("E7").Value = "  -0.13%  "
$cell = $ws.Range("D8")
$cell.Value = "'22.11"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -1.13%  "
$ws.Range("E9").Value = "  -1.35%  "
$ws.Range("E10").Value = "  -0.41%  "
$cell = $ws.Range("D11")
$cell.Value = "'0.0865"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  +0.13%  "
$ws.Range("D12").Value = "1.794.84"
$ws.Range("E12").Value = "  -1.55%  "
$ws.Range("D13").Value = "1.571.84"
$ws.Range("E13").Value = "  -1.81%  "
$ws.Range("E14").Value = "  -1.24%  "
$ws.Range("E15").Value = "  -2.92%  "
$cell = $ws.Range("D16")
$cell.Value = "'63.16"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  -0.41%  "
$ws.Range("D17").Value = "27.445.87"
$ws.Range("E17").Value = "  -0.51%  "
$cell = $ws.Range("D18")
$cell.Value = "'214.04"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -0.73%  "
$cell = $ws.Range("D20")
$cell.Value = "'7.27"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -1.35%  "
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("E22").Value = "  -1.11%  "
$ws.Range("E23").Value = "  +1.38%  "
$cell = $ws.Range("D24")
$cell.Value = "'2.02"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +0.90%  "
$cell = $ws.Range("D25")
$cell.Value = "'153.06"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +0.11%  "
$cell = $ws.Range("D26")
$cell.Value = "'6.80"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +0.52%  "
$ws.Range("E27").Value = "  -0.09%  "
$cell = $ws.Range("D28")
$cell.Value = "'15.06"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -0.35%  "
$ws.Range("E29").Value = "  -1.57%  "
$ws.Range("E30").Value = "  -0.77%  "
$ws.Range("E31").Value = "  +0.83%  "
$ws.Range("E32").Value = "  -1.67%  "
$ws.Range("D33").Value = "1.360.58"
$ws.Range("E33").Value = "  -0.94%  "
$ws.Range("E34").Value = "  -0.44%  "
$ws.Range("E35").Value = "  +1.19%  "
$cell = $ws.Range("D36")
$cell.Value = "'0.971"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +0.27%  "
$ws.Range("E37").Value = "  +0.18%  "
$ws.Range("E38").Value = "  +0.91%  "
$cell = $ws.Range("D39")
$cell.Value = "'0.532"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -1.63%  "
$ws.Range("E40").Value = "  +1.46%  "
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("E42").Value = "  -0.09%  "
$ws.Range("E43").Value = "  -0.15%  "
$cell = $ws.Range("D44")
$cell.Value = "'64.24"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +0.26%  "
$cell = $ws.Range("D45")
$cell.Value = "'5.29"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -1.45%  "
$cell = $ws.Range("D46")
$cell.Value = "'2.17"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -0.09%  "
$ws.Range("D47").Value = "1.708.62"
$ws.Range("E47").Value = "  -1.37%  "
$cell = $ws.Range("D48")
$cell.Value = "'85.33"
$cell.Style = "Normal"
$ws.Range("D49").Value = "0.0₇0998"
$ws.Range("E49").Value = "  -0.74%  "
$cell = $ws.Range("D50")
$cell.Value = "'0.0957"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -1.38%  "
$cell = $ws.Range("D51")
$cell.Value = "'0.0495"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -0.61%  "
